$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N3").Value = 2020
